$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 1760
$ws.Range("C5").Value = 1760

$ws.Range("C6").Select()
